$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Working Papers")
$ws.Activate()

$newAbstract = "Economists often study non-market goods such as health and air quality. This paper introduces a new method to estimate demand for such amenities and applies it to measure the value of a statistical life (VSL) in Kenya. My approach is to update beliefs about the life-saving efficacy of a product (a motorcycle helmet) and elicit product choice. This generates instruments allowing one to use subjective beliefs to estimate demand, rather than assuming rational expectations. This method does not require beliefs to be reported error-free but does require classical mismeasurement. I validate this assumption using features of the experimental design. The estimated VSL is \`$224, near the left tail of Kenyan estimates. Standard methods for estimating VSL produce skewed results, driven by severe violations of rational expectations. These findings help explain low observed demand for many health products and suggest that directing more development aid towards consumption may increase welfare. "

$ws.Range("C3").Value = $newAbstract
$ws.Range("C3").Select() | Out-Null
